$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @(
    "2003","2004","2005","2006","2007","2008","2009","2010","2011","2012",
    "2101","2102","2103","2104","2105","2106","2107","2108","2109","2110","2111","2112",
    "2201","2202","2203","2204","2205","2206","2207","2208","2209","2210","2211","2212",
    "2301","2302","2303","2304","2305","2306","2307","2308","2309","2310","2311","2312",
    "2401"
)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

$ws.Cells.Item(16, 6).Value = 35112
$ws.Cells.Item(62, 6).Value = 10534
